$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-12 18:48:45'
$ws.Range('E3').Value = '2026-02-12 18:48:48'
$ws.Range('O3').Value = '-3.2 °C'
$ws.Range('E4').Value = '2026-02-12 18:48:50'
$ws.Range('J4').Value = '998.7 hPa'
$ws.Range('E5').Value = '2026-02-12 18:48:53'
$ws.Range('E6').Value = '2026-02-12 18:48:55'
$ws.Range('J6').Value = '998.5 hPa'
$ws.Range('O6').Value = '16.1 °C'
$ws.Range('E7').Value = '2026-02-12 18:48:57'
$ws.Range('J7').Value = '1001.4 hPa'
$ws.Range('E8').Value = '2026-02-12 18:49:00'
$ws.Range('J8').Value = '1000.7 hPa'
$ws.Range('E9').Value = '2026-02-12 18:49:02'
$ws.Range('H9').Value = '''69%'
$ws.Range('E10').Value = '2026-02-12 18:49:05'
$ws.Range('E11').Value = '2026-02-12 18:49:07'
$ws.Range('E12').Value = '2026-02-12 18:49:10'
$ws.Range('H12').Value = '''74%'
$ws.Range('O12').Value = '12.7 °C'
$ws.Range('E13').Value = '2026-02-12 18:49:12'
$ws.Range('H13').Value = '''50%'
$ws.Range('J13').Value = '1001.3 hPa'
$ws.Range('O13').Value = '7.8 °C'
$ws.Range('E14').Value = '2026-02-12 18:49:14'
$ws.Range('O14').Value = '17.1 °C'
$ws.Range('E15').Value = '2026-02-12 18:49:17'
$ws.Range('E16').Value = '2026-02-12 18:49:19'
$ws.Range('E17').Value = '2026-02-12 18:49:22'
$ws.Range('E18').Value = '2026-02-12 18:49:24'
$ws.Range('J18').Value = '998.9 hPa'
$ws.Range('O18').Value = '17.0 °C'
$ws.Range('E19').Value = '2026-02-12 18:49:27'
$ws.Range('E20').Value = '2026-02-12 18:49:29'
$ws.Range('H20').Value = '''84%'
$ws.Range('E21').Value = '2026-02-12 18:49:32'
$ws.Range('J21').Value = '1001.7 hPa'
$ws.Range('E22').Value = '2026-02-12 18:49:34'
$ws.Range('E23').Value = '2026-02-12 18:49:37'
$ws.Range('E24').Value = '2026-02-12 18:49:39'
$ws.Range('J24').Value = '1006.3 hPa'
$ws.Range('E25').Value = '2026-02-12 18:49:42'
$ws.Range('E26').Value = '2026-02-12 18:49:44'
$ws.Range('H26').Value = '''45%'
$ws.Range('J26').Value = '997.9 hPa'
$ws.Range('E27').Value = '2026-02-12 18:49:47'
$ws.Range('H27').Value = '''63%'
$ws.Range('E28').Value = '2026-02-12 18:49:49'
$ws.Range('J28').Value = '998.2 hPa'
$ws.Range('E29').Value = '2026-02-12 18:49:52'
$ws.Range('H29').Value = '''55%'
$ws.Range('O29').Value = '15.2 °C'
$ws.Range('E30').Value = '2026-02-12 18:49:54'
$ws.Range('J30').Value = '998.7 hPa'
$ws.Range('E31').Value = '2026-02-12 18:49:57'
$ws.Range('J31').Value = '998.1 hPa'
$ws.Range('E32').Value = '2026-02-12 18:49:59'
$ws.Range('E33').Value = '2026-02-12 18:50:02'
$ws.Range('J33').Value = '1000.9 hPa'
$ws.Range('E34').Value = '2026-02-12 18:50:04'
$ws.Range('O34').Value = '0.4 °C'
$ws.Range('E35').Value = '2026-02-12 18:50:07'
$ws.Range('J35').Value = '1007.7 hPa'
$ws.Range('E36').Value = '2026-02-12 18:50:09'
$ws.Range('H36').Value = '''63%'
$ws.Range('J36').Value = '999.1 hPa'
$ws.Range('E37').Value = '2026-02-12 18:50:12'
$ws.Range('J37').Value = '999.5 hPa'
$ws.Range('O37').Value = '10.4 °C'
$ws.Range('E38').Value = '2026-02-12 18:50:14'
$ws.Range('E39').Value = '2026-02-12 18:50:16'
$ws.Range('E40').Value = '2026-02-12 18:50:19'
$ws.Range('H40').Value = '''54%'
$ws.Range('J40').Value = '1002.4 hPa'
$ws.Range('O40').Value = '10.0 °C'
$ws.Range('E41').Value = '2026-02-12 18:50:21'
$ws.Range('J41').Value = '1005.3 hPa'
$ws.Range('E42').Value = '2026-02-12 18:50:24'
$ws.Range('O42').Value = '14.6 °C'
$ws.Range('E43').Value = '2026-02-12 18:50:26'
$ws.Range('E44').Value = '2026-02-12 18:50:29'
$ws.Range('E45').Value = '2026-02-12 18:50:32'
$ws.Range('J45').Value = '1004.4 hPa'
$ws.Range('E46').Value = '2026-02-12 18:50:34'
$ws.Range('J46').Value = '1007.1 hPa'
